$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume symbol data (scraper refresh).
# Cells already carry numeric-looking text (e.g. "303.07", "2.46%"); keep them as
# TEXT (not auto-converted to number/percent) by forcing the Text number format
# before assigning, matching the workbook author's original inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.46%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.34%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.086"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.93%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07707"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.49%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.425"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.69%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.622"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.13%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "11.92%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1264"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.83%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1859"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.50%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09176"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.16%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04171"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.32%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.13%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.51%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005752"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.05%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1,898.93%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.344"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.42%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3347"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.25%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.652"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "7.94%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1369"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.65%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04150"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.37%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001283"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.40%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004450"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "15.15%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001346"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.16%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02459"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.20%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05283"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.28%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005949"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.40%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007651"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.89%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.27%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007374"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.02%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007558"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.76%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3012"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.50%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006704"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.56%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.28%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.04429"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.62%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.14%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.28%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.28%"
